$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and 1h volume-change (E) columns with freshly scraped values.
# Price cells are plain text in this sheet (e.g. "1.00", "56.700.79") - several of them
# look numeric to Excel and would otherwise be auto-converted to numbers (losing
# significant trailing/leading zeros or the multi-dot thousands formatting), so every
# Price cell is forced to Text format before the value is written and the cell style is
# restored to Normal afterwards so no visible formatting change is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.660.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.347.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "515.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.344.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.05%  "
$ws.Range("E10").Value = "  +6.31%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.339"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.739.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.681.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.332.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("E20").Value = "  +2.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "319.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.49%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.996"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  +5.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "170.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  +9.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0735"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.17%  "
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.943"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.86%  "
$ws.Range("E37").Value = "  +2.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.47%  "
$ws.Range("E39").Value = "  +7.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.72%  "
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "137.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "275.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.558"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.66%  "
$ws.Range("E49").Value = "  +4.15%  "
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.28%  "
